$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "helfen",
    "wundern",
    "tropfen",
    "brauen",
    "lockern",
    "schauen",
    "kümmern",
    "achten",
    "lieben",
    "mauern",
    "erben",
    "reiten",
    "gelten",
    "buchen",
    "krachen",
    "nutzen",
    "heulen",
    "trennen",
    "kichern",
    "kosten",
    "tragen",
    "schaden",
    "küssen",
    "sterben",
    "drehen",
    "bergen",
    "weichen",
    "spannen",
    "quälen",
    "zielen",
    "gnaden",
    "sinken",
    "sprengen",
    "hauen",
    "fischen",
    "foltern",
    "warnen",
    "meinen",
    "schreiten",
    "stillen",
    "liefern",
    "bauen",
    "hören",
    "enden",
    "heben",
    "biegen",
    "regeln",
    "spinnen",
    "dienen",
    "schwächen",
    "regnen",
    "träumen",
    "wählen",
    "füttern",
    "bluten",
    "kehren",
    "räumen",
    "schulden",
    "wenden",
    "schwimmen",
    "werfen",
    "betteln",
    "segeln",
    "jubeln",
    "zünden",
    "stoppen",
    "beten",
    "danken",
    "fließen",
    "streichen",
    "trauen",
    "spielen",
    "herzen",
    "wirken",
    "lohnen",
    "greifen",
    "stehlen",
    "klingen",
    "fassen",
    "scheinen",
    "gründen",
    "bitten",
    "drucken",
    "posten",
    "reizen",
    "flüchten",
    "fallen",
    "saufen",
    "ärgern",
    "arten",
    "handeln",
    "pflanzen",
    "stammen",
    "feiern",
    "formen",
    "platzen",
    "steuern",
    "wärmen",
    "baden",
    "zögern",
    "wollen",
    "stecken",
    "münzen",
    "wüten",
    "folgen",
    "nerven",
    "lenken",
    "decken",
    "loben",
    "führen",
    "atmen",
    "ändern",
    "lügen",
    "schlucken",
    "klappen",
    "schwören",
    "sichern",
    "sammeln",
    "tollen",
    "seufzen",
    "dringen",
    "sperren",
    "tanzen",
    "wechseln",
    "klingeln",
    "äußern",
    "altern",
    "rufen",
    "klettern",
    "irren",
    "schrecken",
    "rauchen",
    "werden",
    "teilen",
    "heilen",
    "zeigen",
    "stechen",
    "leiden",
    "schwingen",
    "schenken",
    "planen",
    "jagen",
    "stärken",
    "sorgen",
    "wahren",
    "schmecken",
    "schütteln",
    "knurren",
    "graben",
    "brauchen",
    "siegen",
    "mühen",
    "leihen",
    "malen",
    "spenden",
    "fangen",
    "löschen",
    "knarren",
    "backen",
    "meistern",
    "flehen",
    "filmen",
    "treiben",
    "töten",
    "öffnen",
    "rasen",
    "lösen",
    "machen",
    "mögen",
    "dürfen",
    "opfern",
    "freuen",
    "pfeifen",
    "fahren",
    "grüßen",
    "wehtun",
    "geben",
    "wecken",
    "spüren",
    "bellen",
    "wachsen",
    "kranken",
    "ehren",
    "liegen",
    "fällen",
    "scheitern",
    "duschen",
    "lesen",
    "suchen",
    "warten",
    "boxen",
    "runden"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

Write-Host "Updated $($words.Length) word cells in column B"
